# Weekly cryptos-list refresh (Price / Volume(1h) columns).
# Mirrors the GitHub Actions commit "Updated cryptos list on Wed Jun 14 23:24:40 UTC 2023".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices such as "1.000" / "14.43" / "25.087.01" as literal
# TEXT (the source site already formats them); temporarily force the whole
# Price column to Text so Excel does not reinterpret these as numbers/dates,
# then restore the original (default) style once every value is written.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "25.087.01"
$ws.Range("E2").Value = "  -3.01%  "
$ws.Range("D3").Value = "1.649.76"
$ws.Range("E3").Value = "  -4.99%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "237.17"
$ws.Range("E5").Value = "  -2.37%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.4790"
$ws.Range("E7").Value = "  -7.97%  "
$ws.Range("D8").Value = "0.2623"
$ws.Range("E8").Value = "  -4.33%  "
$ws.Range("D9").Value = "0.06045"
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("D10").Value = "0.07103"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("D11").Value = "1.659.05"
$ws.Range("E11").Value = "  -4.48%  "
$ws.Range("D12").Value = "14.43"
$ws.Range("E12").Value = "  -3.74%  "
$ws.Range("D13").Value = "0.6179"
$ws.Range("E13").Value = "  -3.89%  "
$ws.Range("D14").Value = "4.560"
$ws.Range("D15").Value = "73.00"
$ws.Range("E15").Value = "  -5.52%  "
$ws.Range("D16").Value = "1.000"
$ws.Range("D17").Value = "0.9990"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "25.048.83"
$ws.Range("E18").Value = "  -3.29%  "
$ws.Range("D19").Value = "11.36"
$ws.Range("E19").Value = "  -3.43%  "
$ws.Range("D20").Value = "0.000006555"
$ws.Range("E20").Value = "  -3.26%  "
$ws.Range("D21").Value = "4.416"
$ws.Range("E21").Value = "  +3.14%  "
$ws.Range("D22").Value = "1.864.85"
$ws.Range("E22").Value = "  -4.96%  "
$ws.Range("D23").Value = "8.468"
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("D24").Value = "5.240"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").Value = "133.46"
$ws.Range("E25").Value = "  -2.67%  "
$ws.Range("D26").Value = "14.70"
$ws.Range("E26").Value = "  -3.29%  "
$ws.Range("D27").Value = "1.400"
$ws.Range("E27").Value = "  -7.40%  "
$ws.Range("D28").Value = "1.687"
$ws.Range("E28").Value = "  -4.86%  "
$ws.Range("D29").Value = "101.76"
$ws.Range("E29").Value = "  -3.13%  "
$ws.Range("D30").Value = "3.788"
$ws.Range("E30").Value = "  -4.09%  "
$ws.Range("D31").Value = "0.07917"
$ws.Range("E31").Value = "  -4.06%  "
$ws.Range("D32").Value = "3.554"
$ws.Range("E32").Value = "  -2.75%  "
$ws.Range("D33").Value = "0.04544"
$ws.Range("E33").Value = "  -2.65%  "
$ws.Range("D34").Value = "2.608"
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("D35").Value = "0.9382"
$ws.Range("E35").Value = "  -5.25%  "
$ws.Range("D36").Value = "0.5769"
$ws.Range("E36").Value = "  -6.91%  "
$ws.Range("D37").Value = "2.620"
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("D38").Value = "0.01535"
$ws.Range("E38").Value = "  -4.08%  "
$ws.Range("D40").Value = "0.8351"
$ws.Range("E40").Value = "  +11.88%  "
$ws.Range("D41").Value = "1.819"
$ws.Range("E41").Value = "  -5.45%  "
$ws.Range("D42").Value = "98.61"
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("D43").Value = "0.3705"
$ws.Range("E43").Value = "  -4.03%  "
$ws.Range("D44").Value = "4.810"
$ws.Range("E44").Value = "  -3.85%  "
$ws.Range("D45").Value = "0.1119"
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("D46").Value = "6.036"
$ws.Range("E46").Value = "  -3.63%  "
$ws.Range("D47").Value = "0.05145"
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("E48").Value = "  -5.35%  "
$ws.Range("D49").Value = "29.42"
$ws.Range("E49").Value = "  -3.88%  "
$ws.Range("E50").Value = "  +0.01%  "
# Row 51: USDD was replaced by Decentraland in the rankings
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "0.3330"
$ws.Range("E51").Value = "  -2.49%  "

# Restore the Price column to its original (default/General) style.
$priceRange.Style = "Normal"

